# Add new data row for patient "mar" (June 2, 2022) at the top of that
# patient's block (row 13), pushing the existing rows 13:25 down to 14:26.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row before row 13 - shifts rows 13:25 down to 14:26,
# carrying their formatting (including the date-number-format on B/F) along.
$ws.Rows("13:13").Insert()

# Populate the newly inserted row 13 with the new observation.
$ws.Range("A13").Value = "mar"
$ws.Range("B13").Value = 44714.007638888892
$ws.Range("C13").Value = "crp"
$ws.Range("D13").Value = 2.6
$ws.Range("E13").Value = 1
$ws.Range("F13").Value = 44710.472916666666
$ws.Range("G13").Value = "uc"

$ws.Range("D14").Select()
